$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: new section header "Autonomous Bot (self-deciding)" ---
# Clone the formatting of row 7 (the existing "Baseline Control (random tools)"
# section header row) onto row 18, so we get the same bold-header / shaded-row
# style family, then re-theme the fill color to a new (blue) theme tint and
# overwrite the cell values for the new block.
$ws.Range("A7:H7").Copy() | Out-Null
$ws.Range("A18:H18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18:H18").Interior.ThemeColor = 3
$ws.Range("A18:H18").Interior.TintAndShade = 0.79998168889431442

$ws.Range("A18").Value = "Autonomous Bot (self-deciding)"
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = ""

# --- Rows 19-20: data rows underneath the new header, same shaded style as
# the rest of the new block (re-use row 18's non-header cell style). ---
$ws.Range("B18:H18").Copy() | Out-Null
$ws.Range("B19:H20").PasteSpecial(-4122) | Out-Null
$ws.Range("A19:A20").Copy() | Out-Null
$ws.Range("A19:A20").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").Copy() | Out-Null
$ws.Range("A19:A20").PasteSpecial(-4122) | Out-Null

$ws.Range("A19").Value = "T1"
$ws.Range("B19").Value = 720
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 0.55000000000000004
$ws.Range("E19").Value = 0.47
$ws.Range("F19").Value = 0.25
$ws.Range("G19").Value = 0.45
$ws.Range("H19").Value = 0.28999999999999998

$ws.Range("A20").Value = "T2"
$ws.Range("B20").Value = 920
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = 0.71
$ws.Range("E20").Value = 0.22
$ws.Range("F20").Value = 0.16
$ws.Range("G20").Value = 0.45
$ws.Range("H20").Value = 0.25

# --- Row 21: just a plain (unstyled) label cell ---
$ws.Range("A21").Value = "T3"

# Match the author's final selection state from the diff.
$ws.Range("B21").Select() | Out-Null
